$d = $word.ActiveDocument

# Locate the "LOM3071: Tratamento de Minerios (Requisito fraco)" paragraph,
# then remove the three paragraphs that follow it:
#   - the blank spacer paragraph
#   - "Ver no Jupiter Salvar em pdf Salvar em docx"
#   - the "(c) 2020 . Contact: ..." footer paragraph
# while leaving the trailing blank paragraph and the page-break paragraph intact.

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "LOM3071: Tratamento de Min") {
        $target = $i
        break
    }
}

if ($target -ne $null) {
    $firstToRemove = $d.Paragraphs.Item($target + 1)
    $lastToRemove = $d.Paragraphs.Item($target + 3)
    $r = $d.Range($firstToRemove.Range.Start, $lastToRemove.Range.End)
    $r.Delete()
}
